{"js": "// Map of old multiplication expressions -> new ones, per the commit diff.\nconst replacements = [\n  [\"95\u00d793=\", \"85\u00d746=\"],\n  [\"46\u00d766=\", \"89\u00d770=\"],\n  [\"24\u00d768=\", \"51\u00d721=\"],\n  [\"28\u00d746=\", \"37\u00d724=\"],\n  [\"53\u00d778=\", \"50\u00d733=\"],\n  [\"79\u00d783=\", \"32\u00d728=\"],\n  [\"80\u00d739=\", \"92\u00d790=\"],\n  [\"90\u00d772=\", \"68\u00d720=\"],\n  [\"85\u00d756=\", \"50\u00d783=\"],\n  [\"38\u00d736=\", \"29\u00d727=\"],\n  [\"30\u00d775=\", \"39\u00d772=\"],\n  [\"92\u00d721=\", \"25\u00d777=\"],\n  [\"79\u00d717=\", \"31\u00d711=\"],\n  [\"64\u00d723=\", \"89\u00d771=\"],\n  [\"49\u00d753=\", \"14\u00d727=\"],\n  [\"87\u00d742=\", \"46\u00d788=\"],\n  [\"89\u00d740=\", \"11\u00d776=\"],\n  [\"51\u00d785=\", \"22\u00d711=\"],\n  [\"14\u00d731=\", \"20\u00d792=\"],\n  [\"55\u00d778=\", \"57\u00d736=\"],\n  [\"48\u00d728=\", \"82\u00d764=\"],\n  [\"19\u00d797=\", \"85\u00d711=\"],\n  [\"35\u00d798=\", \"59\u00d759=\"],\n  [\"89\u00d713=\", \"43\u00d744=\"],\n  [\"35\u00d754=\", \"93\u00d736=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"95\u00d793=\", \"85\u00d746=\"),\n    @(\"46\u00d766=\", \"89\u00d770=\"),\n    @(\"24\u00d768=\", \"51\u00d721=\"),\n    @(\"28\u00d746=\", \"37\u00d724=\"),\n    @(\"53\u00d778=\", \"50\u00d733=\"),\n    @(\"79\u00d783=\", \"32\u00d728=\"),\n    @(\"80\u00d739=\", \"92\u00d790=\"),\n    @(\"90\u00d772=\", \"68\u00d720=\"),\n    @(\"85\u00d756=\", \"50\u00d783=\"),\n    @(\"38\u00d736=\", \"29\u00d727=\"),\n    @(\"30\u00d775=\", \"39\u00d772=\"),\n    @(\"92\u00d721=\", \"25\u00d777=\"),\n    @(\"79\u00d717=\", \"31\u00d711=\"),\n    @(\"64\u00d723=\", \"89\u00d771=\"),\n    @(\"49\u00d753=\", \"14\u00d727=\"),\n    @(\"87\u00d742=\", \"46\u00d788=\"),\n    @(\"89\u00d740=\", \"11\u00d776=\"),\n    @(\"51\u00d785=\", \"22\u00d711=\"),\n    @(\"14\u00d731=\", \"20\u00d792=\"),\n    @(\"55\u00d778=\", \"57\u00d736=\"),\n    @(\"48\u00d728=\", \"82\u00d764=\"),\n    @(\"19\u00d797=\", \"85\u00d711=\"),\n    @(\"35\u00d798=\", \"59\u00d759=\"),\n    @(\"89\u00d713=\", \"43\u00d744=\"),\n    @(\"35\u00d754=\", \"93\u00d736=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
